$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the model numbers (column A) for the new rows first, then quantities,
# then destinations, then prices - matching the order in which new values were
# typed (this is also the order the shared-string table records them in).
$ws.Range("A2").Value() = "H52414130"
$ws.Range("A3").Value() = "T41.1.123.57"
$ws.Range("A4").Value() = "H70545540"
$ws.Range("A5").Value() = "T063.209.16.038.00"
$ws.Range("A6").Value() = "C032.607.11.051.00"

$ws.Range("B2").Value() = 1
$ws.Range("B3").Value() = 1
$ws.Range("B4").Value() = 1
$ws.Range("B5").Value() = 1
$ws.Range("B6").Value() = 1

$ws.Range("C2").Value() = "11/12/2023 POP"
$ws.Range("C3").Value() = "11/12/2023 POP"
$ws.Range("C4").Value() = "11/12/2023 POP"
$ws.Range("C5").Value() = "11/12/2023 POP"
$ws.Range("C6").Value() = "11/12/2023 POP"

$ws.Range("D2").Value() = 3889
$ws.Range("D3").Value() = 1929
$ws.Range("D4").Value() = 4989
$ws.Range("D5").Value() = 1399
$ws.Range("D6").Value() = 4749

# Row 2's helper-average formula now references a deleted cell.
$ws.Range("E2").Formula = "=#REF!/B2"

# Give the new rows the same look as row 2 (model column in blue, qty column
# font) by copying just the formatting down from row 2. Column D also picks
# up row 2's (unstyled) look so the new cells don't inherit the column's
# default style.
$ws.Range("A2").Copy()
$ws.Range("A3:A6").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B3:B6").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D3:D6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C10").Select()
